# Update the dSF (column F) values for specific rows on Sheet1.
# These values were re-pulled from source data, so column F ("dSF")
# no longer mirrors column E ("dS0") for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -1
$ws.Range("F10").Value = 2
$ws.Range("F12").Value = -3
$ws.Range("F17").Value = 1
$ws.Range("F19").Value = -1
$ws.Range("F24").Value = -8
$ws.Range("F26").Value = -6
$ws.Range("F27").Value = -11
$ws.Range("F29").Value = -8
$ws.Range("F30").Value = -5
$ws.Range("F34").Value = -10
$ws.Range("F35").Value = 1
$ws.Range("F36").Value = 9
